# CryCompanywiseStockReport_1 - stock quantity/value corrections
# Applies the per-item quantity (F) and value (G) adjustments, the
# associated rate (D/E) and code (B) corrections for a few rows whose
# data was mis-ordered, and recomputes the affected Sub Total / Grand
# Total rows (column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 3M INDIA LTD ---
$ws.Range("F9").Value  = 27
$ws.Range("G9").Value  = 798.39
$ws.Range("B10").Value = 34053.62        # Sub Total

# --- BLUE STAR LIMITED ---
$ws.Range("F112").Value = 6
$ws.Range("G112").Value = 283.8
$ws.Range("B114").Value = 217709.07      # Sub Total

# --- DAB-Real Activ Coconut Water Tetra 1000ml (rows 163/164 swapped) ---
$ws.Range("B163").Value = 64329
$ws.Range("E163").Value = 128.32
$ws.Range("F163").Value = 2
$ws.Range("G163").Value = 241.38

$ws.Range("B164").Value = 57552
$ws.Range("E164").Value = 136.86
$ws.Range("F164").Value = -5
$ws.Range("G164").Value = -603.45

# --- HUL-Bru Inst Poly 50g (rows 277/278 swapped) ---
$ws.Range("B277").Value = 61610
$ws.Range("E277").Value = 122.71
$ws.Range("F277").Value = -58
$ws.Range("G277").Value = -5957.18

$ws.Range("B278").Value = 63565
$ws.Range("E278").Value = 109.19
$ws.Range("F278").Value = 60
$ws.Range("G278").Value = 6162.6

# --- HUL-Kissan Pineapple Jam 500G (rows 294/295/296 rotated) ---
$ws.Range("B294").Value = 63571
$ws.Range("F294").Value = 0
$ws.Range("G294").Value = 0

$ws.Range("B295").Value = 57802
$ws.Range("E295").Value = 162.71
$ws.Range("F295").Value = -79
$ws.Range("G295").Value = -11334.92

$ws.Range("B296").Value = 63531
$ws.Range("E296").Value = 152.53
$ws.Range("F296").Value = 75
$ws.Range("G296").Value = 10761

# --- HUL-knorr schezwan 200g pch (rows 299/300 swapped) ---
$ws.Range("B299").Value = 55356
$ws.Range("E299").Value = 54.04
$ws.Range("F299").Value = -158
$ws.Range("G299").Value = -7527.12

$ws.Range("B300").Value = 63510
$ws.Range("E300").Value = 50.66
$ws.Range("F300").Value = 130
$ws.Range("G300").Value = 6193.2

# --- HUL-Srf Xl Qw Det Liq 1000Ml ---
$ws.Range("F325").Value = 150
$ws.Range("G325").Value = 20673

# --- HUL-Surf Excel Easywash Pwd 1.5Kg ---
$ws.Range("F329").Value = 12
$ws.Range("G329").Value = 1932

$ws.Range("B339").Value = 228496.26      # Sub Total

# --- JLM-MBD Shiny Toothbrush Safari (rows 356/357 swapped) ---
$ws.Range("B356").Value = 31930
$ws.Range("E356").Value = 26.8
$ws.Range("F356").Value = -62
$ws.Range("G356").Value = -1390.04

$ws.Range("B357").Value = 63681
$ws.Range("E357").Value = 23.84
$ws.Range("F357").Value = 0
$ws.Range("G357").Value = 0

# --- JYOTHY-Pril Tamarind Speckles 425ml ---
$ws.Range("F383").Value = 56
$ws.Range("G383").Value = 4663.12

$ws.Range("B395").Value = 210531.23      # Sub Total

# --- KUS-Floor Wiper (rows 420/421 swapped) ---
$ws.Range("B420").Value = 47097
$ws.Range("D420").Value = 112.28
$ws.Range("E420").Value = 134.16
$ws.Range("F420").Value = 15
$ws.Range("G420").Value = 1684.2

$ws.Range("B421").Value = 58047
$ws.Range("D421").Value = 105.54
$ws.Range("E421").Value = 126.1
$ws.Range("F421").Value = 41
$ws.Range("G421").Value = 4327.14

# --- CHUK-Kashmiri mirch sabut 100gm ---
$ws.Range("F442").Value = 12
$ws.Range("G442").Value = 611.4

$ws.Range("B448").Value = 32499.25       # Sub Total

# --- LOREAL-TR5 sh 340ml Renovation bottle ---
$ws.Range("F459").Value = 7
$ws.Range("G459").Value = 1829.1

$ws.Range("B460").Value = 40295.5        # Sub Total

# --- CRE-Cremica Chocolate Cream 150Gm (rows 472/473 swapped) ---
$ws.Range("B472").Value = 45695
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28

$ws.Range("B473").Value = 64915
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0

# --- RANGA-Cycle Brand Woods 24gms ---
$ws.Range("F498").Value = 35
$ws.Range("G498").Value = 2157.75

$ws.Range("B508").Value = 3311.01        # Sub Total

# --- UNB-McvitieS Cookie Butter (70 Gms) ---
$ws.Range("F549").Value = 61
$ws.Range("G549").Value = 786.29

$ws.Range("B555").Value = 3389.63        # Sub Total

# --- RAJ-Pocha S-42 ---
$ws.Range("F583").Value = 35
$ws.Range("G583").Value = 634.9

$ws.Range("B585").Value = 17284.83       # Sub Total

# --- SRL-20W g-line eco LH 6500k led batten ---
$ws.Range("F673").Value = 0
$ws.Range("G673").Value = 0

$ws.Range("B682").Value = 4353.81        # Sub Total

# --- TATA CONSUMER PRODUCT LIMITED items ---
$ws.Range("F701").Value = 0
$ws.Range("G701").Value = 0

$ws.Range("F702").Value = 29
$ws.Range("G702").Value = 1387.94

$ws.Range("F703").Value = 10
$ws.Range("G703").Value = 815.6

$ws.Range("F704").Value = 48
$ws.Range("G704").Value = 6870.24

$ws.Range("F705").Value = 0
$ws.Range("G705").Value = 0

$ws.Range("F706").Value = 40
$ws.Range("G706").Value = 4071.6

$ws.Range("F707").Value = 0
$ws.Range("G707").Value = 0

$ws.Range("F710").Value = 0
$ws.Range("G710").Value = 0

$ws.Range("F711").Value = 26
$ws.Range("G711").Value = 969.28

$ws.Range("F712").Value = 0
$ws.Range("G712").Value = 0

$ws.Range("F714").Value = 4
$ws.Range("G714").Value = 459.44

$ws.Range("F716").Value = 107
$ws.Range("G716").Value = 14446.07

$ws.Range("F717").Value = 0
$ws.Range("G717").Value = 0

$ws.Range("F718").Value = 18
$ws.Range("G718").Value = 2172.78

$ws.Range("B719").Value = 35358.93       # Sub Total

# --- Tip Top Food Tech India ---
$ws.Range("F721").Value = 4
$ws.Range("G721").Value = 997.92

$ws.Range("F744").Value = 50
$ws.Range("G744").Value = 12090

$ws.Range("B746").Value = 40266.6        # Sub Total

# --- VVD AND SONS PRIVATE LIMITED items ---
$ws.Range("F771").Value = 2416
$ws.Range("G771").Value = 394073.76

$ws.Range("F777").Value = 19
$ws.Range("G777").Value = 2443.02

$ws.Range("B778").Value = 600334.65      # Sub Total

# --- WIPRO ENTERPRISES PVT LTD ---
$ws.Range("F781").Value = 9
$ws.Range("G781").Value = 1314.09

$ws.Range("B795").Value = 53463.64       # Sub Total

# --- Report Grand Total ---
$ws.Range("B796").Value = 2230971.48     # Sub Total (grand)
$ws.Range("B797").Value = 2230971.48     # Grand Total

Write-Host "Applied stock report corrections"
